$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the original B and C column values for rows 17-31 before any writes,
# since rows 16-30 need to be overwritten with values taken from the row below.
$origB = @{}
$origC = @{}
for ($r = 17; $r -le 31; $r++) {
    $origB[$r] = $ws.Cells.Item($r, 2).Value()
    $origC[$r] = $ws.Cells.Item($r, 3).Value()
}

# Rows 16-22: only column B is replaced with the next row's original B value.
for ($r = 16; $r -le 22; $r++) {
    $ws.Cells.Item($r, 2).Value = $origB[$r + 1]
}

# Rows 23-30: both column B and column C are replaced with the next row's original values.
for ($r = 23; $r -le 30; $r++) {
    $ws.Cells.Item($r, 2).Value = $origB[$r + 1]
    $ws.Cells.Item($r, 3).Value = $origC[$r + 1]
}

# Row 31 (the old last row) is removed entirely now that its contents have migrated to row 30.
$ws.Rows.Item(31).Delete()
